# Feature/Add JWT Authentication and User login/Registration Endpoints
# Reorganize the Features sheet: regroup Identity, split DB Deployment /
# Resources Management into Local/Cloud variants, add a Pipeline row at
# the end, and highlight rows backed by Azure resources.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row order / content (column A) -------------------------------
$ws.Range("A1").Value  = "Main Feature"
$ws.Range("A2").Value  = "Create Repo"
$ws.Range("A3").Value  = "Create Front-End"
$ws.Range("A4").Value  = "Identity"
$ws.Range("A5").Value  = "Create Back-End"
$ws.Range("A6").Value  = "Create Model Project"
$ws.Range("A7").Value  = "Create DB Project"
$ws.Range("A8").Value  = "Webhost (Azure)"
$ws.Range("A9").Value  = "Resources Management [Blob Storage Local]"
$ws.Range("A10").Value = "DB Deployment [Local]"
$ws.Range("A11").Value = "Resources Management [Blob Storage Cloud]"
$ws.Range("A12").Value = "DB Deployment [Cloud]"
$ws.Range("A13").Value = "Pipeline"

# --- Highlight the rows associated with the Azure-hosted resources ----
# (xlThemeColorAccent6 = 10 -> theme index 9, the workbook's green accent)
$ws.Range("B2").Interior.ThemeColor = 10
$ws.Range("B4").Interior.ThemeColor = 10
$ws.Range("B5").Interior.ThemeColor = 10
$ws.Range("B6").Interior.ThemeColor = 10
$ws.Range("B7").Interior.ThemeColor = 10
$ws.Range("B10").Interior.ThemeColor = 10

# --- Column sizing / view state ----------------------------------------
# Widen column A so the new, longer labels (e.g. the "Resources
# Management [Blob Storage ...]" rows) remain fully visible, mirroring
# the best-fit width Excel computed for the longest label.
$ws.Columns.Item(1).ColumnWidth = 37.3

[void]$ws.Range("F10").Select()
